try {
    $wb = $excel.ActiveWorkbook

    # The backward-elimination summaries were regenerated (statsmodels OLS
    # output) at a later run; only the embedded "run at" Date/Time stamp
    # text needs to be refreshed to match, everything else (coefficients,
    # p-values, etc.) stays as-is.
    $oldDate = "Sun, 29 Dec 2019"
    $newDate = "Wed, 01 Jan 2020"
    $oldTime = "16:11:38"
    $newTime = "23:19:14"

    for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
        $ws = $wb.Worksheets.Item($i)
        $used = $ws.UsedRange
        foreach ($cell in $used.Cells) {
            $text = $cell.Value2
            if ($text -ne $null -and ($text -is [string]) -and $text -like "*$oldDate*") {
                $text = $text -replace [regex]::Escape($oldDate), $newDate
                $text = $text -replace [regex]::Escape($oldTime), $newTime
                $cell.Value2 = $text
            }
        }
    }
}
catch {
    Write-Output "Error updating workbook: $_"
}
